$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 160, shifting rows 160:278 down to 161:279
$ws.Rows.Item(160).Insert()

# Populate the new row 160 with data (constant columns copied from neighboring rows,
# plus the new data point's own values)
$ws.Range("A160").Value = 8
$ws.Range("B160").Value = "Terminal La Palmera de La Serena"
$ws.Range("C160").Value = "Coquimbo"
$ws.Range("D160").Value = 44942
$ws.Range("E160").Value = 4
$ws.Range("F160").Value = 100112037
$ws.Range("G160").Value = "Cebollín"
$ws.Range("H160").Value = "Sin especificar"
$ws.Range("I160").Value = "Primera"
$ws.Range("J160").Value = 1400
$ws.Range("K160").Value = 1200
$ws.Range("L160").Value = 1400
$ws.Range("M160").Value = 1300
$ws.Range("N160").Value = "`$/paquete 6 unidades"
$ws.Range("O160").Value = "Provincia del Elquí"
$ws.Range("P160").Value = 217
$ws.Range("Q160").Value = 6
$ws.Range("R160").Value = "Hortaliza"

# Match the date cell style used by the other date cells in column D (style index 2)
$ws.Range("D160").NumberFormat = $ws.Range("D161").NumberFormat
